# Update the NATMI ligand-receptor table with newly computed TPM-based values.
# The underlying raw inputs that changed are the ligand (F8) expression values
# for the "ECs" sending cluster and the receptor (Ldlr) expression values for
# the "ECs" target cluster. All dependent specificity / edge-weight columns
# are recomputed from those raw values, mirroring the original analysis script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New raw ligand expression values (F8) for rows where Sending cluster = "ECs"
$ligandAvgECs = 0.190922
$ligandTotECs = 0.572766

# New raw receptor expression values (Ldlr) for rows where Target cluster = "ECs"
$receptorAvgECs = 1.090291
$receptorTotECs = 3.270873

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $sendingCluster = $ws.Cells.Item($r, 1).Value2
    $targetCluster  = $ws.Cells.Item($r, 4).Value2

    if ($sendingCluster -eq "ECs") {
        $ws.Cells.Item($r, 7).Value2 = $ligandAvgECs   # G - Ligand average expression value
        $ws.Cells.Item($r, 8).Value2 = $ligandTotECs   # H - Ligand total expression value
    }

    if ($targetCluster -eq "ECs") {
        $ws.Cells.Item($r, 13).Value2 = $receptorAvgECs  # M - Receptor average expression value
        $ws.Cells.Item($r, 14).Value2 = $receptorTotECs  # N - Receptor total expression value
    }
}

# Recompute specificity (derived from the three sending/target clusters) and
# edge-weight columns, which depend on the raw values updated above.
$gByCluster = @{}
$mByCluster = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $sendingCluster = $ws.Cells.Item($r, 1).Value2
    $targetCluster  = $ws.Cells.Item($r, 4).Value2
    $gByCluster[$sendingCluster] = $ws.Cells.Item($r, 7).Value2
    $mByCluster[$targetCluster]  = $ws.Cells.Item($r, 13).Value2
}

$gTotal = 0.0
foreach ($v in $gByCluster.Values) { $gTotal += $v }
$mTotal = 0.0
foreach ($v in $mByCluster.Values) { $mTotal += $v }

for ($r = 2; $r -le $lastRow; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    $h = $ws.Cells.Item($r, 8).Value2
    $m = $ws.Cells.Item($r, 13).Value2
    $n = $ws.Cells.Item($r, 14).Value2

    $i = $g / $gTotal
    $j = $i
    $o = $m / $mTotal
    $p = $o

    $ws.Cells.Item($r, 9).Value2  = $i   # I - Ligand derived specificity of average expression value
    $ws.Cells.Item($r, 10).Value2 = $j   # J - Ligand derived specificity of total expression value
    $ws.Cells.Item($r, 15).Value2 = $o   # O - Receptor derived specificity of average expression value
    $ws.Cells.Item($r, 16).Value2 = $p   # P - Receptor derived specificity of total expression value

    $ws.Cells.Item($r, 17).Value2 = $g * $m   # Q - Edge average expression weight
    $ws.Cells.Item($r, 18).Value2 = $h * $n   # R - Edge total expression weight
    $ws.Cells.Item($r, 19).Value2 = $i * $o   # S - Edge average expression derived specificity
    $ws.Cells.Item($r, 20).Value2 = $j * $p   # T - Edge total expression derived specificity
}
